$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-03-03 Monday" "2025-03-04 Tuesday"

Replace-Text "605×4=" "943×9="
Replace-Text "815×2=" "828×8="
Replace-Text "905×9=" "684×4="
Replace-Text "931×8=" "195×3="
Replace-Text "715×9=" "212×4="
Replace-Text "257×5=" "522×6="
Replace-Text "584×6=" "362×6="
Replace-Text "746×9=" "535×9="
Replace-Text "188×8=" "590×8="
Replace-Text "106×8=" "527×5="
Replace-Text "514×2=" "477×9="
Replace-Text "712×7=" "213×9="
Replace-Text "716×6=" "237×9="
Replace-Text "264×8=" "711×6="
Replace-Text "571×3=" "899×3="
Replace-Text "961×7=" "227×2="
Replace-Text "800×3=" "352×4="
Replace-Text "566×2=" "321×9="
Replace-Text "372×4=" "953×3="
Replace-Text "659×2=" "964×9="
Replace-Text "332×7=" "240×5="
Replace-Text "858×2=" "287×3="
Replace-Text "927×5=" "843×6="
Replace-Text "956×7=" "696×3="
Replace-Text "151×7=" "699×7="
